$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B values
$ws.Range("B1").Value = -9.8782200000000007
$ws.Range("B2").Value = -9.8804700000000008
$ws.Range("B3").Value = -9.9316899999999997
$ws.Range("B4").Value = -9.7557399999999994
$ws.Range("B5").Value = -10.022270000000001
$ws.Range("B6").Value = -10.30518
$ws.Range("B7").Value = -10.43038
$ws.Range("B8").Value = -10.455249999999999
$ws.Range("B9").Value = -10.54266
$ws.Range("B10").Value = -10.64167
$ws.Range("B11").Value = -10.572710000000001

# Update the selection to whole columns A:B
$ws.Range("A1:B1048576").Select()
